$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fmt = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- Extend the thin bottom-border row (row 2) and header row (row 3) with
# --- three new columns (AK:AM), cloning formatting from the last existing
# --- triplet (AH:AJ) so the new cells pick up the right style indices.
$ws.Range("AH2").Copy()
$ws.Range("AK2").PasteSpecial($fmt)
$ws.Range("AH2").Copy()
$ws.Range("AL2").PasteSpecial($fmt)
$ws.Range("AH2").Copy()
$ws.Range("AM2").PasteSpecial($fmt)

$ws.Range("AH3").Copy()
$ws.Range("AK3").PasteSpecial($fmt)
$ws.Range("AH3").Copy()
$ws.Range("AL3").PasteSpecial($fmt)
$ws.Range("AH3").Copy()
$ws.Range("AM3").PasteSpecial($fmt)

# --- Clone formatting into the new AK:AM data cells for every data row too,
# --- so the new columns inherit the same per-row style as AH:AJ (s="1" for
# --- rows 5-11, s="5" for the bottom total row 12).
foreach ($r in 5..12) {
    $ws.Range("AH$r").Copy()
    $ws.Range("AK$r").PasteSpecial($fmt)
    $ws.Range("AH$r").Copy()
    $ws.Range("AL$r").PasteSpecial($fmt)
    $ws.Range("AH$r").Copy()
    $ws.Range("AM$r").PasteSpecial($fmt)
}

$excel.CutCopyMode = $false

# --- Header row 3: years 2018-2021 relabeled / shifted one slot to the right
$ws.Range("AC3").Value = "2018 факт "
$ws.Range("AF3").Value = "2019 факт"
$ws.Range("AH3").Value = "2020 утв."
$ws.Range("AI3").Value = "2020 факт"
$ws.Range("AJ3").Value = "откл. от утв., %"
$ws.Range("AK3").Value = "2021 утв."
$ws.Range("AL3").Value = "2021 факт"
$ws.Range("AM3").Value = "откл. от утв., %"

# --- Row 5 data (2019 план/факт/откл. + new 2020/2021 triplets)
$ws.Range("AE5").Value = 43737.8
$ws.Range("AF5").Value = 43258.3
$ws.Range("AG5").Value = 98.9
$ws.Range("AH5").Value = 46293.5
$ws.Range("AI5").Value = 47153.5
$ws.Range("AJ5").Value = 101.9
$ws.Range("AK5").Value = 47483.3
$ws.Range("AL5").Value = 52020.5
$ws.Range("AM5").Value = 109.6

# --- Row 6
$ws.Range("AE6").Value = 6265.4
$ws.Range("AF6").Value = 4434.6000000000004
$ws.Range("AG6").Value = 70.8
$ws.Range("AH6").Value = 7935.8
$ws.Range("AI6").Value = 3895.8
$ws.Range("AJ6").Value = 49.1
$ws.Range("AK6").Value = 8997
$ws.Range("AL6").Value = 6212.4
$ws.Range("AM6").Value = 69

# --- Row 7
$ws.Range("AE7").Value = 728.5
$ws.Range("AF7").Value = 695.7
$ws.Range("AG7").Value = 95.5
$ws.Range("AH7").Value = 746.9
$ws.Range("AI7").Value = 583.20000000000005
$ws.Range("AJ7").Value = 78.099999999999994
$ws.Range("AK7").Value = 639.20000000000005
$ws.Range("AL7").Value = 600.79999999999995
$ws.Range("AM7").Value = 94

# --- Row 8
$ws.Range("AE8").Value = 1249
$ws.Range("AF8").Value = 1244.7
$ws.Range("AG8").Value = 99.7
$ws.Range("AH8").Value = 1249
$ws.Range("AI8").Value = 1207.5999999999999
$ws.Range("AJ8").Value = 96.7
$ws.Range("AK8").Value = 1208.0999999999999
$ws.Range("AL8").Value = 1332.7
$ws.Range("AM8").Value = 110.3

# --- Row 9
$ws.Range("AE9").Value = 2582.6
$ws.Range("AF9").Value = 2477.5
$ws.Range("AG9").Value = 95.9
$ws.Range("AH9").Value = 3109
$ws.Range("AI9").Value = 3225.2
$ws.Range("AJ9").Value = 103.7
$ws.Range("AK9").Value = 3131.3
$ws.Range("AL9").Value = 4833.7
$ws.Range("AM9").Value = 154.4

# --- Row 10
$ws.Range("AE10").Value = 2686.4
$ws.Range("AF10").Value = 2829
$ws.Range("AG10").Value = 105.3
$ws.Range("AH10").Value = 2993.4
$ws.Range("AI10").Value = 2624.5
$ws.Range("AJ10").Value = 87.7
$ws.Range("AK10").Value = 2798.4
$ws.Range("AL10").Value = 3088
$ws.Range("AM10").Value = 110.3

# --- Row 11
$ws.Range("AE11").Value = 23397.4
$ws.Range("AF11").Value = 24364.799999999999
$ws.Range("AG11").Value = 104.1
$ws.Range("AH11").Value = 30085.9
$ws.Range("AI11").Value = 29223.5
$ws.Range("AJ11").Value = 97.1
$ws.Range("AK11").Value = 30439.7
$ws.Range("AL11").Value = 30705.3
$ws.Range("AM11").Value = 100.9

# --- Row 12
$ws.Range("AE12").Value = 13137.1
$ws.Range("AF12").Value = 10924.7
$ws.Range("AG12").Value = 83.2
$ws.Range("AH12").Value = 12158.7
$ws.Range("AI12").Value = 10980.3
$ws.Range("AJ12").Value = 90.3
$ws.Range("AK12").Value = 11664.9
$ws.Range("AL12").Value = 11939.1
$ws.Range("AM12").Value = 102.4

# --- Mirror the author's active-cell selection recorded in the saved view
$ws.Range("AF4").Select()
